$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.959.89'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.826.67'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.82'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4646'
$ws.Range("E7").Value = '  -1.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3702'
$ws.Range("E8").Value = '  +1.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07369'
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8729'
$ws.Range("E10").Value = '  -0.81%  '
$ws.Range("E11").Value = '  -1.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07832'
$ws.Range("E12").Value = '  +6.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.806.09'
$ws.Range("E13").Value = '  -5.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.581'
$ws.Range("E14").Value = '  +0.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.345'
$ws.Range("E15").Value = '  -0.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.61'
$ws.Range("E16").Value = '  -1.73%  '
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008819'
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.009'
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.61'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.910.63'
$ws.Range("E21").Value = '  -2.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.149'
$ws.Range("E22").Value = '  -1.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.57'
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.009.91'
$ws.Range("E24").Value = '  -3.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.61'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.840'
$ws.Range("E26").Value = '  -2.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.29'
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.076'
$ws.Range("E28").Value = '  -2.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.112'
$ws.Range("E29").Value = '  -1.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.21'
$ws.Range("E30").Value = '  -0.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08867'
$ws.Range("E31").Value = '  -0.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.970'
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7289'
$ws.Range("E33").Value = '  -1.76%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.445'
$ws.Range("E34").Value = '  -1.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.134'
$ws.Range("E35").Value = '  -2.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.479'
$ws.Range("E36").Value = '  +1.41%  '
$ws.Range("E37").Value = '  -1.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01947'
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("E39").Value = '  -1.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.922'
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.196'
$ws.Range("E41").Value = '  -1.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5185'
$ws.Range("E42").Value = '  -1.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8620'
$ws.Range("E43").Value = '  -14.52%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1628'
$ws.Range("E44").Value = '  -0.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.209'
$ws.Range("E45").Value = '  -2.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4832'
$ws.Range("E46").Value = '  -0.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.010'
$ws.Range("E47").Value = '  +0.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.18'
$ws.Range("E48").Value = '  -1.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.45'
$ws.Range("E49").Value = '  -1.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.626'
$ws.Range("E50").Value = '  -1.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06224'
$ws.Range("E51").Value = '  -0.81%  '
